# BUG: Don't extract header names if none specified (#23703)
# Add a new worksheet "index_col_none" used by the pandas Excel-reader
# tests, demonstrating a MultiIndex-column header (two header rows) that
# is read back with index_col=None.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it lands at the
# end of the tab strip.
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$newSheet = $wb.Sheets.Add($null, $lastSheet)
$newSheet.Name = "index_col_none"

# Two header rows: top level "A"/"A"/"B"/"B", second level "key"/"val"/"key"/"val"
$newSheet.Range("A1").Value = "A"
$newSheet.Range("B1").Value = "A"
$newSheet.Range("C1").Value = "B"
$newSheet.Range("D1").Value = "B"

$newSheet.Range("A2").Value = "key"
$newSheet.Range("B2").Value = "val"
$newSheet.Range("C2").Value = "key"
$newSheet.Range("D2").Value = "val"

$headerRange = $newSheet.Range("A1:D2")
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.Font.Bold = $true

# Two data rows
$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = 2
$newSheet.Cells.Item(3, 3).Value = 3
$newSheet.Cells.Item(3, 4).Value = 4

$newSheet.Cells.Item(4, 1).Value = 1
$newSheet.Cells.Item(4, 2).Value = 2
$newSheet.Cells.Item(4, 3).Value = 3
$newSheet.Cells.Item(4, 4).Value = 4

$newSheet.Range("A3:D4").HorizontalAlignment = -4108   # xlCenter

$newSheet.PageSetup.Orientation = 1   # xlPortrait

# Make the new sheet the active one, with the same lingering selection
# the authored workbook shows.
$wb.ActiveSheet = $newSheet
$null = $newSheet.Range("G23").Select()
